$d = $word.ActiveDocument
$r = $d.Content
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>

    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="44"/>
        </w:rPr>
        <w:t>Embracing Diversity in a Globalized World: Navigating Challenges and Celebrating Unity</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="36"/>
        </w:rPr>
        <w:t>Hannah Davies</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>hannah</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>davies@eduworld</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>org</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>In a world where borders are increasingly blurred and cultures intertwine, embracing diversity has become a cornerstone of global harmony</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> As citizens of a globalized society, it is imperative that we understand, appreciate, and celebrate the rich tapestry of cultures, beliefs, and perspectives that make our world a vibrant and dynamic place</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> In this essay, we will delve into the multifaceted nature of diversity, explore the challenges it presents, and highlight the immense benefits that stem from fostering inclusivity and understanding</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>The beauty of diversity lies in its multifaceted nature</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> It encompasses differences in race, ethnicity, gender, sexual orientation, religion, language, and socioeconomic status, among others</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Each individual brings a unique set of experiences, perspectives, and talents to the table, creating a rich and dynamic social fabric</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Diversity challenges us to step outside our comfort zones, question our assumptions, and engage with those who hold different beliefs and values</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> By doing so, we broaden our horizons, foster empathy, and gain a deeper understanding of the human experience</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>However, embracing diversity is not without its challenges</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Misunderstandings, prejudice, and discrimination can arise when people from different backgrounds interact</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> These challenges can be daunting, but they also present opportunities for growth and learning</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> By promoting open dialogue, encouraging inclusivity, and challenging stereotypes, we can create a society where diversity is celebrated and everyone feels valued and respected</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>Introduction Continued:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>The immense benefits that stem from fostering inclusivity and understanding are undeniable</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> A diverse society is a more vibrant, innovative, and resilient one</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> When people from different backgrounds come together, they bring fresh ideas, perspectives, and solutions to complex problems</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> This cross-pollination of ideas leads to groundbreaking discoveries, transformative technologies, and creative expressions that benefit all of humanity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Furthermore, a diverse society is a more tolerant and peaceful </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>one</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> When people understand and respect each other's differences, they are less likely to resort to violence or conflict</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Instead, they are more likely to work together to build a better future for themselves and for generations to come</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>Introduction Continued:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>The journey towards a truly diverse and inclusive society is an ongoing one, but it is a journey worth taking</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> By challenging our biases, educating ourselves about different cultures, and actively promoting inclusivity, we can create a world where everyone feels valued, respected, and empowered</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> A world where diversity is not just tolerated but celebrated, and where the unique contributions of each individual are recognized and appreciated</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> In this world, the challenges of diversity will be outweighed by the boundless opportunities it presents, leading to a more harmonious, equitable, and prosperous global community</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>Summary</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>In conclusion, embracing diversity in a globalized world presents both challenges and immense benefits</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> By fostering inclusivity, understanding, and respect, we can create a society where everyone feels valued and empowered</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> The challenges of diversity, such as misunderstandings and prejudice, can be overcome through open dialogue, education, and a commitment to building a more just and equitable world</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> The beauty of diversity lies in its multifaceted nature, and it is through celebrating our differences and coming together as a global community that we can unlock the full potential of our shared humanity</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
$r.InsertXML($xml)
